# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
# on Sheet1, matching the latest crypto price/volume snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.276.20"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.858.08"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.18"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07298"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8921"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.11"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07876"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.856.45"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.414"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.521"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.70"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008937"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "27.307.59"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "2.076.73"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.034"
$ws.Range("E25").Value = "  +9.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.65"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.06"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.050"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08842"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.146"
$ws.Range("E32").Value = "  +6.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7718"
$ws.Range("E33").Value = "  +5.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.170"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.709"
$ws.Range("E36").Value = "  +9.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.106"
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01944"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.946"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.069"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5119"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.521"
$ws.Range("E44").Value = "  +4.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4801"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.36"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.97"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.647"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06205"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.59"
$ws.Range("E51").Value = "  +1.38%  "
